$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

function Set-EmptyText($cell) {
    # Forces the cell to be stored as an empty Text-typed cell (matching the
    # source file's empty inlineStr cells) instead of leaving it unset, while
    # avoiding the quotePrefix formatting side effect of the leading "'" trick.
    $cell.Value = "'"
    $cell.Style = "Normal"
}

function Set-LiteralText($cell, [string]$text) {
    # Forces a literal text value (e.g. a date-looking string) to be stored
    # as text rather than being parsed into a date/number, again without
    # leaving a quotePrefix style behind.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$ws.Cells.Item($row, 1).Value = 112439580
$ws.Cells.Item($row, 2).Value = 56575
$ws.Cells.Item($row, 3).Value = "Ovaliderad"
$ws.Cells.Item($row, 4).Value = "NT"
$ws.Cells.Item($row, 5).Value = 103021
$ws.Cells.Item($row, 6).Value = "Talltita"
$ws.Cells.Item($row, 7).Value = "Poecile montanus"
$ws.Cells.Item($row, 8).Value = "(Conrad von Baldenstein, 1827)"
Set-EmptyText $ws.Cells.Item($row, 9)
Set-EmptyText $ws.Cells.Item($row, 11)
Set-EmptyText $ws.Cells.Item($row, 12)
$ws.Cells.Item($row, 13).Value = "spel/sång"
Set-EmptyText $ws.Cells.Item($row, 14)
$ws.Cells.Item($row, 16).Value = "Stigsbo, Dlr"
$ws.Cells.Item($row, 17).Value = 570818
$ws.Cells.Item($row, 18).Value = 6702190
$ws.Cells.Item($row, 19).Value = 15
$ws.Cells.Item($row, 20).Value = "Dalarna"
$ws.Cells.Item($row, 21).Value = "Hedemora"
$ws.Cells.Item($row, 22).Value = "Dalarna"
$ws.Cells.Item($row, 23).Value = "Husby"
Set-LiteralText $ws.Cells.Item($row, 25) "2023-10-01"
Set-LiteralText $ws.Cells.Item($row, 27) "2023-10-01"
$ws.Cells.Item($row, 30).Value = $false
$ws.Cells.Item($row, 31).Value = $false
$ws.Cells.Item($row, 33).Value = $false
Set-EmptyText $ws.Cells.Item($row, 46)
$ws.Cells.Item($row, 49).Value = "Philipp Weiss"
$ws.Cells.Item($row, 50).Value = "Philipp Weiss"
Set-EmptyText $ws.Cells.Item($row, 51)
